$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before the current row 23 (old rows 23-31 shift down to 26-34)
$ws.Rows("23:25").Insert()

# --- New row 23 ---
$ws.Range("A23").Value = 5
$ws.Range("B23").Value = "Macroferia Regional de Talca"
$ws.Range("C23").Value = "Maule"
$ws.Range("D23").Value = "12/23/2021"
$ws.Range("E23").Value = 7
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100103
$ws.Range("H23").Value = "Frutos de hueso (carozo)"
$ws.Range("I23").Value = 100103003
$ws.Range("J23").Value = "Damasco"
$ws.Range("K23").Value = "Castle Brite"
$ws.Range("L23").Value = "Especial"
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 13000
$ws.Range("O23").Value = 13000
$ws.Range("P23").Value = 13000
$ws.Range("Q23").Value = "$/caja 10 kilos"
$ws.Range("R23").Value = "Región de O'Higgins"
$ws.Range("S23").Value = 1300
$ws.Range("T23").Value = 10

# --- New row 24 ---
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "Macroferia Regional de Talca"
$ws.Range("C24").Value = "Maule"
$ws.Range("D24").Value = "12/23/2021"
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100103
$ws.Range("H24").Value = "Frutos de hueso (carozo)"
$ws.Range("I24").Value = 100103003
$ws.Range("J24").Value = "Damasco"
$ws.Range("K24").Value = "Castle Brite"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 120
$ws.Range("N24").Value = 11000
$ws.Range("O24").Value = 11000
$ws.Range("P24").Value = 11000
$ws.Range("Q24").Value = "$/caja 10 kilos"
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 1100
$ws.Range("T24").Value = 10

# --- New row 25 ---
$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Macroferia Regional de Talca"
$ws.Range("C25").Value = "Maule"
$ws.Range("D25").Value = "12/23/2021"
$ws.Range("E25").Value = 7
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100103
$ws.Range("H25").Value = "Frutos de hueso (carozo)"
$ws.Range("I25").Value = 100103003
$ws.Range("J25").Value = "Damasco"
$ws.Range("K25").Value = "Castle Brite"
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 150
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 14000
$ws.Range("P25").Value = 14000
$ws.Range("Q25").Value = "$/caja 15 kilos"
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 933
$ws.Range("T25").Value = 15
